$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.847493588924408
$ws.Range("B1").Value = 2.282288074493408
$ws.Range("C1").Value = 3.405540704727173
$ws.Range("D1").Value = 1.661390066146851
$ws.Range("E1").Value = 1.206437826156616
